# Commit: "update on minerals csv"
# - Rename the worksheet from "Sheet1" to "minerals"
# - Scroll the sheet view back to the top (A1) instead of A17,
#   while keeping the existing selection on F30

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to match the data it now represents
$ws.Name = "minerals"

# Make sure it's the active sheet and reset the visible top-left
# cell back to A1 (removes the stale topLeftCell="A17" scroll position)
$ws.Activate()
$ws.Range("A1").Select()

# Restore the original active selection on F30
$ws.Range("F30").Select()
